$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B18: was stored as text "4", should be a real number 4
$ws.Range("B18").Value = 4

# Add new row 19 of annotation data
$ws.Range("A19").Value = "Sunsi Wu"

# B19 keeps "3" as text (not a number), matching the source data
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "3"

$ws.Range("C19").Value = "无"
$ws.Range("D19").Value = "SMY"
$ws.Range("E19").Value = "EXP"
$ws.Range("F19").Value = "2bb8b329-99fa-4c06-a5b4-7897e3cce401"
$ws.Range("G19").Value = "S1PWi_lC-_annotated.xlsx"
$ws.Range("H19").Value = "Each network is trained with 50 epochs."
